# Fixed naive component forecaster bug - Presentation state 11.02.
#
# A new "naive" error value is inserted into column B for every data row
# (rows 2..16). All previously existing values in that row shift one
# column to the right (B->C, C->D, ...). Column K is the right-hand edge
# of the table, so for the fully populated rows (2..6) the old value that
# was in column K simply falls off the end and is dropped; the other,
# shorter rows just grow by one column to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values to insert into column B for each row.
$newValues = @{
    2  = -0.6603092772102132
    3  = -0.15162438770796
    4  = -0.2053460154962278
    5  = 0.6162032393936197
    6  = 1.652643173475852
    7  = 0.3110387314724781
    8  = 0.2388379152847414
    9  = 0.6508000635779043
    10 = 0.2387740594105157
    11 = 0.3465902496671606
    12 = 0.00230005330798793
    13 = -0.1902738424076751
    14 = -0.3325070745318338
    15 = 0.1656141382254278
    16 = -0.09587373626955231
}

$firstDataCol = 2   # column B
$lastTableCol = 11  # column K (the table never grows past this column)

for ($row = 2; $row -le 16; $row++) {

    # Find the last currently-populated column in B:K for this row.
    $lastUsedCol = $firstDataCol - 1
    for ($c = $firstDataCol; $c -le $lastTableCol; $c++) {
        $cellVal = $ws.Cells.Item($row, $c).Value()
        if ($cellVal -ne $null -and $cellVal -ne "") {
            $lastUsedCol = $c
        }
    }

    if ($lastUsedCol -ge $firstDataCol) {
        # Capture the existing values before overwriting anything.
        $oldValues = @()
        for ($c = $firstDataCol; $c -le $lastUsedCol; $c++) {
            $oldValues += , $ws.Cells.Item($row, $c).Value()
        }

        # Write them back one column to the right, dropping whatever
        # would fall past the last table column (K).
        for ($i = 0; $i -lt $oldValues.Count; $i++) {
            $destCol = $firstDataCol + $i + 1
            if ($destCol -le $lastTableCol) {
                $ws.Cells.Item($row, $destCol).Value = $oldValues[$i]
            }
        }
    }

    # Insert the new value into column B.
    $ws.Cells.Item($row, $firstDataCol).Value = $newValues[$row]
}
